$d = $word.ActiveDocument

# Locate the run containing "Projeto Recicla " (with trailing space) that
# follows "Recursos:" and split it so the period is inserted right after
# "Projeto Recicla" and before the existing trailing space.
$found = $d.Content.Find.Execute("Projeto Recicla ", $true, $false, $false, $false, $false, `
                                  $true, 1, $false, "Projeto Recicla. ", 2)
